$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 218.44444
$ws.Range("I33").Value = 113.27273
$ws.Range("J33").Value = 681.2
$ws.Range("K33").Value = 113.27273
$ws.Range("L33").Value = 681.2
$ws.Range("M33").Value = 115.72727
$ws.Range("N33").Value = -1139.2

$ws.Range("H40").Value = 1961.826
$ws.Range("I40").Value = 1874
$ws.Range("J40").Value = 2029.3846
$ws.Range("K40").Value = 1874
$ws.Range("L40").Value = 2029.3846
$ws.Range("M40").Value = -1699
$ws.Range("N40").Value = -2379.3846

$ws.Range("H96").Value = 1052.2106
$ws.Range("I96").Value = 616.13336
$ws.Range("J96").Value = 2687.5
$ws.Range("K96").Value = 1848.40008
$ws.Range("L96").Value = 8062.5
$ws.Range("M96").Value = -475.4000800000001

$ws.Range("H129").Value = 1985.4166
$ws.Range("I129").Value = 1711.3636
$ws.Range("J129").Value = 5000
$ws.Range("K129").Value = 5134.0908
$ws.Range("L129").Value = 15000
$ws.Range("M129").Value = -134.0907999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17889.328
$ws.Range("I32").Value = 18757.902
$ws.Range("J32").Value = 228.33333
$ws.Range("K32").Value = 18757.902
$ws.Range("L32").Value = 228.33333
$ws.Range("M32").Value = -18470.902

$ws.Range("H61").Value = 6668114
$ws.Range("I61").Value = 8334642.5
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 8334642.5
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -8334430.5

$ws.Range("H102").Value = 33523.91
$ws.Range("I102").Value = 36781.3
$ws.Range("J102").Value = 950
$ws.Range("K102").Value = 36781.3
$ws.Range("L102").Value = 950
$ws.Range("M102").Value = -35159.3

$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H122").Value = 3556.4443
$ws.Range("I122").Value = 3438.5625
$ws.Range("J122").Value = 4499.5
$ws.Range("K122").Value = 10315.6875
$ws.Range("L122").Value = 13498.5
$ws.Range("M122").Value = -7865.6875
$ws.Range("N122").Value = -18398.5

$ws.Range("H136").Value = 6668114
$ws.Range("I136").Value = 8334642.5
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 25003927.5
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -25001377.5

$ws.Range("H139").Value = 77405.336
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 77405.336
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 77405.336
$ws.Range("N139").Value = -87685.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 374.875
$ws.Range("I22").Value = 374.875
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 374.875
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -201.875

$ws.Range("H134").Value = 748647.5
$ws.Range("I134").Value = 681947.0600000001
$ws.Range("J134").Value = 1526818.6
$ws.Range("K134").Value = 2045841.18
$ws.Range("L134").Value = 4580455.800000001
$ws.Range("M134").Value = -2043306.18

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 995
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 995
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 995
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -1221

$ws.Range("H7").Value = 329.2353
$ws.Range("I7").Value = 351.1111
$ws.Range("J7").Value = 304.625
$ws.Range("K7").Value = 351.1111
$ws.Range("L7").Value = 304.625
$ws.Range("M7").Value = -238.1111
$ws.Range("N7").Value = -530.625

$ws.Range("H86").Value = 140357.2
$ws.Range("I86").Value = 6813.4287
$ws.Range("J86").Value = 257208
$ws.Range("K86").Value = 6813.4287
$ws.Range("L86").Value = 257208
$ws.Range("M86").Value = -5690.4287

$ws.Range("H89").Value = 140357.2
$ws.Range("I89").Value = 6813.4287
$ws.Range("J89").Value = 257208
$ws.Range("K89").Value = 34067.14350000001
$ws.Range("L89").Value = 1286040
$ws.Range("M89").Value = -28451.14350000001

$ws.Range("H95").Value = 19999.5
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 19999.5
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 19999.5
$ws.Range("N95").Value = -25491.5

$ws.Range("H103").Value = 10222
$ws.Range("I103").Value = 10222
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 10222
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -9050

$ws.Range("H105").Value = 87869
$ws.Range("I105").Value = 87869
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 87869
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -86122

$ws.Range("H106").Value = 80000
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 80000
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 80000
$ws.Range("N106").Value = -82524

$ws.Range("H107").Value = 743.6061
$ws.Range("I107").Value = 548.7143
$ws.Range("J107").Value = 1084.6666
$ws.Range("K107").Value = 548.7143
$ws.Range("L107").Value = 1084.6666
$ws.Range("M107").Value = 1371.2857

$ws.Range("H108").Value = 80464
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 80464
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 80464
$ws.Range("N108").Value = -88144

$ws.Range("H109").Value = 18000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 18000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 18000
$ws.Range("N109").Value = -20080

$ws.Range("H132").Value = 3844.45
$ws.Range("I132").Value = 2319.5
$ws.Range("J132").Value = 7402.6665
$ws.Range("K132").Value = 6958.5
$ws.Range("L132").Value = 22207.9995
$ws.Range("M132").Value = -4428.5
$ws.Range("N132").Value = -27267.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 15009.167
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 15009.167
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 45027.501
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -46525.501

$ws.Range("H66").Value = 15009.167
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 15009.167
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 135082.503
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -142570.503

$ws.Range("H94").Value = 3874.2856
$ws.Range("I94").Value = 2024
$ws.Range("J94").Value = 8500
$ws.Range("K94").Value = 6072
$ws.Range("L94").Value = 25500
$ws.Range("M94").Value = -5396
$ws.Range("N94").Value = -26852

$ws.Range("H101").Value = 12997.5
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 12997.5
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 38992.5
$ws.Range("N101").Value = -43860.5

$ws.Range("H138").Value = 4495
$ws.Range("I138").Value = 4694
$ws.Range("J138").Value = 3500
$ws.Range("K138").Value = 14082
$ws.Range("L138").Value = 10500
$ws.Range("M138").Value = -8942
$ws.Range("N138").Value = -20780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 9112.125
$ws.Range("I57").Value = 3271.4285
$ws.Range("J57").Value = 49997
$ws.Range("K57").Value = 3271.4285
$ws.Range("L57").Value = 49997
$ws.Range("M57").Value = -2451.4285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4064.1428
$ws.Range("I7").Value = 3827.0908
$ws.Range("J7").Value = 4933.3335
$ws.Range("K7").Value = 3827.0908
$ws.Range("L7").Value = 4933.3335
$ws.Range("M7").Value = -3715.0908
$ws.Range("N7").Value = -5157.3335

$ws.Range("H14").Value = 500
$ws.Range("I14").Value = 500
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 500
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -328

$ws.Range("H22").Value = 1364.8334
$ws.Range("I22").Value = 797.25
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 797.25
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -502.25
$ws.Range("N22").Value = -3090

$ws.Range("H27").Value = 1364.8334
$ws.Range("I27").Value = 797.25
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 797.25
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -690.25
$ws.Range("N27").Value = -2714

$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()

$ws.Range("H46").Value = 6933.75
$ws.Range("I46").Value = 10715.429
$ws.Range("J46").Value = 1639.4
$ws.Range("K46").Value = 10715.429
$ws.Range("L46").Value = 1639.4
$ws.Range("M46").Value = -10527.429
$ws.Range("N46").Value = -2015.4

$ws.Range("H100").Value = 15125.375
$ws.Range("I100").Value = 2999.8333
$ws.Range("J100").Value = 51502
$ws.Range("K100").Value = 2999.8333
$ws.Range("L100").Value = 51502
$ws.Range("M100").Value = -2458.8333

$ws.Range("H122").Value = 4899.452
$ws.Range("I122").Value = 4196.2334
$ws.Range("J122").Value = 6657.5
$ws.Range("K122").Value = 12588.7002
$ws.Range("L122").Value = 19972.5
$ws.Range("M122").Value = -10138.7002
$ws.Range("N122").Value = -24872.5

$ws.Range("H124").Value = 80000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 80000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 80000
$ws.Range("N124").Value = -89820

$ws.Range("H126").Value = 4064.1428
$ws.Range("I126").Value = 3827.0908
$ws.Range("J126").Value = 4933.3335
$ws.Range("K126").Value = 11481.2724
$ws.Range("L126").Value = 14800.0005
$ws.Range("M126").Value = -9011.2724
$ws.Range("N126").Value = -19740.0005
